$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Handback report refresh ---
# Status text: handback just completed and the target is in sync with en-US
# (the "Status" column value is shared by the Overview summary sheet too)
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime refreshed for this handback run
$wsZhCn.Range("K2").Value = "2016-09-02 04:52:54"
$wsDeDe.Range("K2").Value = "2016-09-02 04:53:02"

# Error Detail cleared now that the handback succeeded (no more stale-handback warning)
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Column widths resize to fit the new, longer "Status" text / shorter "Error Detail" text
$wsOverview.Columns("E").ColumnWidth = 29.144371396019366
$wsOverview.Columns("F").ColumnWidth = 29.144371396019366
$wsZhCn.Columns("C").ColumnWidth = 29.144371396019366
$wsZhCn.Columns("P").ColumnWidth = 12.913719540550566
$wsDeDe.Columns("C").ColumnWidth = 29.144371396019366
$wsDeDe.Columns("P").ColumnWidth = 12.913719540550566
